$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new rows 50 and 51 inherit the same style as the rest of column A (bold, centered, bordered)
$ws.Range("A49").Copy()
$ws.Range("A50:A51").PasteSpecial(-4122)

$rows = @(
  @(2, 13, 0.17142, 0.171428956652766, 0.000008956652765984563, 0.999991043427455),
  @(3, 12, 0.15877, 0.1588348147326499, 0.00006481473264993709, 0.9999351894680274),
  @(4, 26, 0.22481, 0.2231413111608817, 0.001668688839118332, 0.9983340910445625),
  @(5, 10, 0.16776, 0.1659385432086418, 0.0018214567913582, 0.9981818548814158),
  @(6, 18, 0.16776, 0.1659385432086418, 0.0018214567913582, 0.9981818548814158),
  @(7, 29, 0.15737, 0.1594464066356196, 0.002076406635619615, 0.9979278958950936),
  @(8, 4, 0.16439, 0.1664767301669309, 0.002086730166930922, 0.9979176152082332),
  @(9, 35, 0.22368, 0.2213441903847347, 0.002335809615265305, 0.9976696336767996),
  @(10, 25, 0.15315, 0.156088591248787, 0.002938591248786954, 0.9970700187684193),
  @(11, 16, 0.15315, 0.156088591248787, 0.002938591248786954, 0.9970700187684193),
  @(12, 44, 0.15315, 0.156088591248787, 0.002938591248786954, 0.9970700187684193),
  @(13, 20, 0.15315, 0.156088591248787, 0.002938591248786954, 0.9970700187684193),
  @(14, 33, 0.16748, 0.1642570642932437, 0.00322293570675633, 0.9967874182377162),
  @(15, 32, 0.16186, 0.1651145987488062, 0.003254598748806181, 0.9967559593019907),
  @(16, 19, 0.17732, 0.1738521128339571, 0.003467887166042932, 0.9965440975138359),
  @(17, 30, 0.1939, 0.1974294051802059, 0.003529405180205891, 0.9964830077105991),
  @(18, 41, 0.18462, 0.1810281619637109, 0.003591838036289097, 0.9964210170906559),
  @(19, 24, 0.18266, 0.178891989034667, 0.003768010965333024, 0.9962461336442577),
  @(20, 28, 0.16018, 0.1641109707039884, 0.003930970703988401, 0.9960844213210877),
  @(21, 15, 0.16861, 0.1726037236736787, 0.003993723673678701, 0.9960221627092793),
  @(22, 40, 0.18378, 0.179775799818457, 0.004004200181542988, 0.9960117694917821),
  @(23, 21, 0.16861, 0.1727121012152924, 0.00410210121529242, 0.9959146572740685),
  @(24, 23, 0.17732, 0.173158783592591, 0.004161216407409046, 0.9958560275587056),
  @(25, 36, 0.20711, 0.2113230416875116, 0.004213041687511632, 0.9958046335661684),
  @(26, 14, 0.16158, 0.1660880669640323, 0.004508066964032331, 0.9955121644989301),
  @(27, 8, 0.16158, 0.1661076915302122, 0.004527691530212236, 0.9954927160611022),
  @(28, 9, 0.17535, 0.1707552234816436, 0.00459477651835638, 0.9954262388917843),
  @(29, 42, 0.17535, 0.1706713675790318, 0.004678632420968221, 0.9953431552438872),
  @(30, 5, 0.22762, 0.2227701717392069, 0.004849828260793138, 0.9951735790519193),
  @(31, 43, 0.20682, 0.2116749473021086, 0.004854947302108586, 0.9951685093304825),
  @(32, 39, 0.15315, 0.1582000249363662, 0.005050024936366221, 0.9949753496730812),
  @(33, 7, 0.17704, 0.1715426682470617, 0.005497331752938261, 0.9945327236787844),
  @(34, 11, 0.17985, 0.1856202330115029, 0.005770233011502895, 0.9942628715564334),
  @(35, 6, 0.21638, 0.2223350944354925, 0.005955094435492525, 0.9940801587780275),
  @(36, 48, 0.17142, 0.1649602682526871, 0.00645973174731293, 0.9935817285644422),
  @(37, 37, 0.17029, 0.1636660381672836, 0.006623961832716407, 0.9934196263114417),
  @(38, 27, 0.17338, 0.1665379255371373, 0.006842074462862752, 0.9932044213919915),
  @(39, 22, 0.1672, 0.1742576594568286, 0.007057659456828597, 0.9929918020179349),
  @(40, 0, 0.19109, 0.1836756537298511, 0.007414346270148892, 0.9926402216748256),
  @(41, 46, 0.20345, 0.2110188873507921, 0.007568887350792058, 0.9924879703553641),
  @(42, 49, 0.21244, 0.2036841098322762, 0.00875589016772374, 0.9913201099958207),
  @(43, 47, 0.18153, 0.1698938860758195, 0.01163611392418051, 0.988497727825232),
  @(44, 38, 0.18153, 0.1696123615841421, 0.01191763841585788, 0.9882227189610858),
  @(45, 31, 0.22902, 0.2138717721090488, 0.01514822789095124, 0.9850778167416765),
  @(46, 17, 0.15456, 0.1356372788601496, 0.01892272113985036, 0.9814286984211307),
  @(47, 2, 0.18856, 0.1659385432086418, 0.02262145679135821, 0.9778789535060834),
  @(48, 45, 0.18125, 0.156975332970231, 0.02427466702976899, 0.9763006273501214),
  @(49, 1, 0.22368, 0.1899030866869787, 0.0337769133130213, 0.9673266902384443),
  @(50, 3, 0.11606, 0.153914430608477, 0.03785443060847701, 0.963526261976563),
  @(51, 34, 0.1377, 0.08446691207348882, 0.05323308792651117, 0.9494574481786262)
)

foreach ($r in $rows) {
  $ws.Cells.Item($r[0], 1).Value = $r[1]
  $ws.Cells.Item($r[0], 2).Value = $r[2]
  $ws.Cells.Item($r[0], 3).Value = $r[3]
  $ws.Cells.Item($r[0], 4).Value = $r[4]
  $ws.Cells.Item($r[0], 5).Value = $r[5]
}
